# Remove the "waste heat" / "recovered heat" recovery connection that ran
# from electricity (power) into the CO2 capture units. This drops both
# rows describing that flow on the "connections" sheet, so no waste heat
# is modeled anywhere in the factory any more.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("connections")

# Rows 25 and 26 are:
#   power -> simple_power -> waste heat -> recovered heat -> CO2 Capture
#   power -> simple_power -> waste heat -> recovered heat -> CO2 Capture Aux
# Deleting the pair shifts every following row up by two.
$ws.Rows("25:26").Delete()

# Leave the sheet focused near where the edit happened, as a user would
# after deleting the rows and continuing to work on the sheet.
$ws.Activate()
$ws.Range("B32").Select()
